# Kevin time log update - add 4 new log entries (rows 28-31)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (number formats / styles) from the last existing data
# row (27) down across the four new rows so the new cells inherit the
# same date / duration formats as the rest of the log.
$ws.Range("A27:E27").Copy() | Out-Null
$ws.Range("A28:E31").PasteSpecial(-4122) | Out-Null

# Row 28 - TUI Implementation
$ws.Range("A28").Value = 44620.708333333336
$ws.Range("B28").Value = 44620.770833333336
$ws.Range("C28").Formula = "=B28-A28"
$ws.Range("D28").Formula = "=C28+D27"
$ws.Range("E28").Value = "TUI Implementation"

# Row 29 - Practice with sub-windows in Ncurses
$ws.Range("A29").Value = 44621.416666666664
$ws.Range("B29").Value = 44621.447916666664
$ws.Range("C29").Formula = "=B29-A29"
$ws.Range("D29").Formula = "=C29+D28"
$ws.Range("E29").Value = "Practice with sub-windows in Ncurses"

# Row 30 - Client meeting
$ws.Range("A30").Value = 44622.416666666664
$ws.Range("B30").Value = 44622.427083333336
$ws.Range("C30").Formula = "=B30-A30"
$ws.Range("D30").Formula = "=C30+D29"
$ws.Range("E30").Value = "Client meeting"

# Row 31 - Team Meeting
$ws.Range("A31").Value = 44622.583333333336
$ws.Range("B31").Value = 44622.59375
$ws.Range("C31").Formula = "=B31-A31"
$ws.Range("D31").Formula = "=C31+D30"
$ws.Range("E31").Value = "Team Meeting"

# Match the author's final selection state
$ws.Range("F30").Select() | Out-Null
